# Update column G ("K" = Strikeouts) values in Sheet1 to reflect
# regenerated save_data (Strike# -> K recalculation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 2
